# Scheduled market-data refresh: update currentAveragePrice / Leve profit figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 40758.305
$ws.Range("J17").Value = 42556.41
$ws.Range("L17").Value = 127669.23
$ws.Range("N17").Value = -128005.23

$ws.Range("H19").Value = 1424.6
$ws.Range("I19").Value = 797
$ws.Range("K19").Value = 797
$ws.Range("M19").Value = -622

$ws.Range("H112").Value = 2347.4
$ws.Range("I112").Value = 1071.75
$ws.Range("J112").Value = 7450
$ws.Range("K112").Value = 3215.25
$ws.Range("L112").Value = 22350
$ws.Range("M112").Value = -2107.25
$ws.Range("N112").Value = -24566

$ws.Range("H132").Value = 1217.6111
$ws.Range("I132").Value = 815.14
$ws.Range("K132").Value = 2445.42
$ws.Range("M132").Value = 84.57999999999993

$ws.Range("H137").Value = 2702.3713
$ws.Range("I137").Value = 2356.4348
$ws.Range("K137").Value = 7069.3044
$ws.Range("M137").Value = -4519.3044

$ws.Range("H138").Value = 2260.65
$ws.Range("I138").Value = 1233.05
$ws.Range("J138").Value = 3288.25
$ws.Range("K138").Value = 3699.15
$ws.Range("L138").Value = 9864.75
$ws.Range("M138").Value = 1440.85
$ws.Range("N138").Value = -20144.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3685.9333
$ws.Range("I2").Value = 3429.6
$ws.Range("K2").Value = 3429.6
$ws.Range("M2").Value = -3316.6

$ws.Range("H32").Value = 42608.79
$ws.Range("I32").Value = 23879.09
$ws.Range("K32").Value = 23879.09
$ws.Range("M32").Value = -23592.09

$ws.Range("H45").Value = 442216.22
$ws.Range("I45").Value = 633854.4
$ws.Range("J45").Value = 4186.143
$ws.Range("K45").Value = 633854.4
$ws.Range("L45").Value = 4186.143
$ws.Range("M45").Value = -633477.4
$ws.Range("N45").Value = -4940.143

$ws.Range("H61").Value = 1315.5238
$ws.Range("I61").Value = 1201.5278
$ws.Range("K61").Value = 1201.5278
$ws.Range("M61").Value = -989.5278000000001

$ws.Range("H74").Value = 1554.8298
$ws.Range("I74").Value = 1438.4878
$ws.Range("K74").Value = 1438.4878
$ws.Range("M74").Value = -564.4878000000001

$ws.Range("H77").Value = 1554.8298
$ws.Range("I77").Value = 1438.4878
$ws.Range("K77").Value = 7192.439
$ws.Range("M77").Value = -2824.439

$ws.Range("H116").Value = 3685.9333
$ws.Range("I116").Value = 3429.6
$ws.Range("K116").Value = 3429.6
$ws.Range("M116").Value = -1135.6

$ws.Range("H122").Value = 2454.2222
$ws.Range("I122").Value = 2254.8462
$ws.Range("J122").Value = 2972.6
$ws.Range("K122").Value = 6764.5386
$ws.Range("L122").Value = 8917.799999999999
$ws.Range("M122").Value = -4314.5386
$ws.Range("N122").Value = -13817.8

$ws.Range("H132").Value = 1320.5106
$ws.Range("I132").Value = 810.5814
$ws.Range("K132").Value = 2431.7442
$ws.Range("M132").Value = 98.25579999999991

$ws.Range("H136").Value = 1315.5238
$ws.Range("I136").Value = 1201.5278
$ws.Range("K136").Value = 3604.5834
$ws.Range("M136").Value = -1054.5834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3685.9333
$ws.Range("I3").Value = 3429.6
$ws.Range("K3").Value = 3429.6
$ws.Range("M3").Value = -3315.6

$ws.Range("H76").Value = 29999.75
$ws.Range("J76").Value = 29999.75
$ws.Range("L76").Value = 29999.75
$ws.Range("N76").Value = -30629.75

$ws.Range("H79").Value = 29999.75
$ws.Range("J79").Value = 29999.75
$ws.Range("L79").Value = 29999.75
$ws.Range("N79").Value = -32183.75

$ws.Range("H99").Value = 1799.5
$ws.Range("I99").Value = 959.4
$ws.Range("K99").Value = 959.4
$ws.Range("M99").Value = 538.6

$ws.Range("H134").Value = 1352.15
$ws.Range("I134").Value = 1352.15
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4056.45
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1521.45
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 624.5
$ws.Range("I22").Value = 593.8889
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 593.8889
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -243.8889
$ws.Range("N22").Value = -1600

$ws.Range("H31").Value = 3568.1333
$ws.Range("I31").Value = 1958
$ws.Range("J31").Value = 10008.667
$ws.Range("K31").Value = 1958
$ws.Range("L31").Value = 10008.667
$ws.Range("M31").Value = -1663
$ws.Range("N31").Value = -10598.667

$ws.Range("H34").Value = 3568.1333
$ws.Range("I34").Value = 1958
$ws.Range("J34").Value = 10008.667
$ws.Range("K34").Value = 1958
$ws.Range("L34").Value = 10008.667
$ws.Range("M34").Value = -1756
$ws.Range("N34").Value = -10412.667

$ws.Range("H58").Value = 1176.6428
$ws.Range("I58").Value = 1009.9
$ws.Range("J58").Value = 1593.5
$ws.Range("K58").Value = 1009.9
$ws.Range("L58").Value = 1593.5
$ws.Range("M58").Value = -806.9
$ws.Range("N58").Value = -1999.5

$ws.Range("H132").Value = 2004.5385
$ws.Range("I132").Value = 2005.4706
$ws.Range("J132").Value = 1998.2
$ws.Range("K132").Value = 6016.4118
$ws.Range("L132").Value = 5994.6
$ws.Range("M132").Value = -3486.4118
$ws.Range("N132").Value = -11054.6

$ws.Range("H134").Value = 1936.3529
$ws.Range("I134").Value = 1874.6
$ws.Range("J134").Value = 2399.5
$ws.Range("K134").Value = 5623.799999999999
$ws.Range("L134").Value = 7198.5
$ws.Range("M134").Value = -3088.799999999999
$ws.Range("N134").Value = -12268.5

$ws.Range("H136").Value = 1176.6428
$ws.Range("I136").Value = 1009.9
$ws.Range("J136").Value = 1593.5
$ws.Range("K136").Value = 3029.7
$ws.Range("L136").Value = 4780.5
$ws.Range("M136").Value = -479.6999999999998
$ws.Range("N136").Value = -9880.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 192.4375
$ws.Range("J61").Value = 313.2857
$ws.Range("L61").Value = 939.8571000000001
$ws.Range("N61").Value = -1369.8571

$ws.Range("H122").Value = 1416.25
$ws.Range("J122").Value = 942.2
$ws.Range("L122").Value = 8479.800000000001
$ws.Range("N122").Value = -13379.8

$ws.Range("H137").Value = 10004049
$ws.Range("I137").Value = 16668209
$ws.Range("J137").Value = 7808.25
$ws.Range("K137").Value = 50004627
$ws.Range("L137").Value = 23424.75
$ws.Range("M137").Value = -49999527
$ws.Range("N137").Value = -33624.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1914.3667
$ws.Range("I102").Value = 1894
$ws.Range("K102").Value = 1894
$ws.Range("M102").Value = -272

$ws.Range("H132").Value = 3889.3635
$ws.Range("I132").Value = 3471.125
$ws.Range("J132").Value = 5004.6665
$ws.Range("K132").Value = 10413.375
$ws.Range("L132").Value = 15013.9995
$ws.Range("M132").Value = -7883.375
$ws.Range("N132").Value = -20073.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 37187.25
$ws.Range("I7").Value = 47083
$ws.Range("K7").Value = 47083
$ws.Range("M7").Value = -46971

$ws.Range("H40").Value = 2804.36
$ws.Range("I40").Value = 2216.842
$ws.Range("K40").Value = 2216.842
$ws.Range("M40").Value = -2080.842

$ws.Range("H46").Value = 2573.6843
$ws.Range("I46").Value = 1915.5834
$ws.Range("K46").Value = 1915.5834
$ws.Range("M46").Value = -1727.5834

$ws.Range("H126").Value = 37187.25
$ws.Range("I126").Value = 47083
$ws.Range("K126").Value = 141249
$ws.Range("M126").Value = -138779

$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040

$ws.Range("H132").Value = 3184.5107
$ws.Range("I132").Value = 2261.8125
$ws.Range("J132").Value = 5152.933
$ws.Range("K132").Value = 6785.4375
$ws.Range("L132").Value = 15458.799
$ws.Range("M132").Value = -4255.4375
$ws.Range("N132").Value = -20518.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 6320387
$ws.Range("I3").Value = 7216156.5
$ws.Range("J3").Value = 50000
$ws.Range("K3").Value = 7216156.5
$ws.Range("L3").Value = 50000
$ws.Range("M3").Value = -7216042.5
$ws.Range("N3").Value = -50228

$ws.Range("H81").Value = 4488.609
$ws.Range("J81").Value = 4123
$ws.Range("L81").Value = 8246
$ws.Range("N81").Value = -10368

$ws.Range("H84").Value = 4488.609
$ws.Range("J84").Value = 4123
$ws.Range("L84").Value = 41230
$ws.Range("N84").Value = -51838

$ws.Range("I107").Value = 1949.5
$ws.Range("J107").Value = 250001010
$ws.Range("K107").Value = 5848.5
$ws.Range("L107").Value = 750003030
$ws.Range("M107").Value = -3928.5
$ws.Range("N107").Value = -750006870

$ws.Range("H132").Value = 6713.6665
$ws.Range("I132").Value = 6713.6665
$ws.Range("K132").Value = 20140.9995
$ws.Range("M132").Value = -17610.9995

$ws.Range("H136").Value = 599.75
$ws.Range("I136").Value = 599.75
$ws.Range("K136").Value = 1799.25
$ws.Range("M136").Value = 750.75
